$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.291.10'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '1.708.44'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''223.61'
$ws.Range("E5").Value = '  -2.37%  '
$ws.Range("D6").Value = '''0.5291'
$ws.Range("E6").Value = '  -2.10%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -4.09%  '
$ws.Range("D9").Value = '''0.06536'
$ws.Range("E9").Value = '  -3.33%  '
$ws.Range("D10").Value = '''20.91'
$ws.Range("E10").Value = '  -2.50%  '
$ws.Range("D11").Value = '''0.07635'
$ws.Range("E11").Value = '  -2.10%  '
$ws.Range("E12").Value = '  -2.94%  '
$ws.Range("D13").Value = '1.707.72'
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("D14").Value = '1.944.66'
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("E15").Value = '  -3.94%  '
$ws.Range("D16").Value = '0.0₅8190'
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").Value = '''67.23'
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("D18").Value = '27.286.09'
$ws.Range("E18").Value = '  -0.64%  '
$ws.Range("D19").Value = '''215.39'
$ws.Range("E19").Value = '  +2.73%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").Value = '''4.665'
$ws.Range("E21").Value = '  -2.69%  '
$ws.Range("D22").Value = '''10.45'
$ws.Range("E22").Value = '  -4.02%  '
$ws.Range("D23").Value = '''5.963'
$ws.Range("E23").Value = '  -3.97%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = '''142.36'
$ws.Range("E25").Value = '  -2.60%  '
$ws.Range("D26").Value = '''1.753'
$ws.Range("E26").Value = '  +8.14%  '
$ws.Range("D27").Value = '''0.1217'
$ws.Range("E27").Value = '  -2.45%  '
$ws.Range("D28").Value = '''7.263'
$ws.Range("E28").Value = '  -2.21%  '
$ws.Range("D29").Value = '''16.30'
$ws.Range("E29").Value = '  -3.01%  '
$ws.Range("D30").Value = '''0.05363'
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("D32").Value = '''3.486'
$ws.Range("E32").Value = '  -4.91%  '
$ws.Range("D33").Value = '''3.416'
$ws.Range("E33").Value = '  -2.95%  '
$ws.Range("D34").Value = '''1.638'
$ws.Range("E34").Value = '  +1.07%  '
$ws.Range("D35").Value = '''2.873'
$ws.Range("E35").Value = '  +0.72%  '
$ws.Range("D36").Value = '''0.9497'
$ws.Range("E36").Value = '  -2.50%  '
$ws.Range("D37").Value = '''2.418'
$ws.Range("E37").Value = '  -1.05%  '
$ws.Range("D38").Value = '''0.5856'
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = '''0.01624'
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").Value = '1.041.58'
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").Value = '''0.8394'
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").Value = '''101.04'
$ws.Range("E44").Value = '  -1.37%  '
$ws.Range("D45").Value = '1.852.21'
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("E46").Value = '  +3.73%  '
$ws.Range("D47").Value = '''58.02'
$ws.Range("E47").Value = '  -2.45%  '
$ws.Range("D48").Value = '''0.4497'
$ws.Range("E48").Value = '  +1.33%  '
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").Value = '''8.082'
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("D51").Value = '''0.05236'
$ws.Range("E51").Value = '  -0.70%  '
